# Updated Masterdata as per 2nd may Data Refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user")

# Data updates
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Update the visible/scrolled view & active selection
# (mirrors the author re-scrolling to row 13 and re-selecting C19 before saving)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
